$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each cell update from the source diff.
# Cells whose new text is a plain numeric literal get their
# NumberFormat set to Text ("@") first so Excel keeps storing
# the original formatted string instead of silently converting
# it into a floating point number (which would also drop
# formatting such as trailing zeros or thousand separators).
$ws.Range("D2").Value = "70.105.99"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "3.615.10"
$ws.Range("E3").Value = "  +3.38%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.01"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "195.11"
$ws.Range("E6").Value = "  -1.31%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -1.85%  "
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.81"
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.51"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").Value = "4.192.65"
$ws.Range("E14").Value = "  +3.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.21"
$ws.Range("E15").Value = "  +4.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "591.98"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.26"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "70.304.76"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("D19").Value = "3.602.07"
$ws.Range("E19").Value = "  +2.88%  "
$ws.Range("E20").Value = "  +1.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.995"
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.88"
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.17"
$ws.Range("E23").Value = "  +3.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "102.90"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.64"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.06"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.86"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.59"
$ws.Range("E28").Value = "  -1.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.87"
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("E30").Value = "  -2.01%  "
$ws.Range("E31").Value = "  -2.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.31"
$ws.Range("E32").Value = "  -3.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.116"
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.19"
$ws.Range("E34").Value = "  -0.55%  "
$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").Value = "0.0₃0872"
$ws.Range("E35").Value = "  +8.60%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "3.931.87"
$ws.Range("E36").Value = "  +6.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "524.55"
$ws.Range("E38").Value = "  +2.72%  "
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.08"
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.392"
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.55"
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.141"
$ws.Range("E45").Value = "  +1.04%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.85"
$ws.Range("E46").Value = "  +0.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.32"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000248"
$ws.Range("E50").Value = "  +3.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.30"
$ws.Range("E51").Value = "  +2.56%  "
